# Applies the "Added Kevin's test case, updated excel file" commit.
#
# Summary of content changes being applied to Sheet1:
#  - Header I1: "Can automate" -> "Testable"
#  - Row 2 (Elgin / Search function loading speed): I2 -> "Yes"
#  - Row 3 (Kevin / Find Library location): H3 -> "Pass", I3 -> "Yes"
#  - Row 4 (Tricia / Subscribe to eNewsletter): G4 -> "Test passes, Successfully
#    subscribes to the newsletter", I4 -> "Yes"
#  - Row 5 (Jun Jie / Submit Feedback): I5 "Yes?" -> "Yes"
#  - Rows 7-11 (Filtering facility type.. through Feedback form validation) are hidden
#  - Row 12: new test case "Confirm Registration with Singpass" filled in, H12 removed
#  - Row 13: new test case "Page efficiency" filled in, H13 removed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header -----------------------------------------------------------
$ws.Range("I1").Value = "Testable"

# --- Row 2 -------------------------------------------------------------
$ws.Range("I2").Value = "Yes"

# --- Row 3 -------------------------------------------------------------
$ws.Range("H3").Value = "Pass"
$ws.Range("I3").Value = "Yes"

# --- Row 4 -------------------------------------------------------------
$ws.Range("G4").Value = "Test passes, Successfully subscribes to the newsletter"
$ws.Range("I4").Value = "Yes"

# --- Row 5 -------------------------------------------------------------
$ws.Range("I5").Value = "Yes"

# --- Hide rows 7 through 11 --------------------------------------------
$ws.Rows("7:11").Hidden = $true

# --- Row 12: new "Confirm Registration with Singpass" test case --------
$ws.Range("C12").Value = "Confirm Registration" + [char]10 + "with Singpass"
$ws.Range("D12").Value = "Test if the registration works for singpass"
$ws.Range("F12").Value = "1. Click on the person icon at the top right of the website." + [char]10 + "2. Click on Apply now!" + [char]10 + "3. Use Singpass " + [char]10 + "4. Create an account with singpass"
$ws.Range("G12").Value = "Registration successful"
$ws.Range("I12").Value = "No"
$ws.Range("H12").Style = "Normal"
$ws.Range("H12").ClearContents()

# --- Row 13: new "Page efficiency" test case ----------------------------
$ws.Range("C13").Value = "Page efficiency"
$ws.Range("D13").Value = "Testing the efficiency to drill down from main page all the way down"
$ws.Range("I13").Value = "No"
$ws.Range("H13").Style = "Normal"
$ws.Range("H13").ClearContents()
